# "Funcionalidade de alterar informação" - update a despesa value and
# register two new (mostly blank) rows at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NEW / RS / BR row: DESPESAS (E3) corrected from 750000 to 50
$ws.Range("E3").Value = 50

# New row 6: only DESPESAS (E6) filled in
$ws.Range("E6").Value = 50

# New row 7: only FATURAMENTO (D7) filled in
$ws.Range("D7").Value = 0
